$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Sun Apr 16 19:06:08 UTC 2023 with GitHub Actions
# Map of row -> (Price, Volume(1h)) ; $null means "leave unchanged"
$updates = @{
    2  = @("30.698.35", "  +0.90%  ")
    3  = @("2.144.74", "  +1.84%  ")
    4  = @("1.009", "  +0.46%  ")
    5  = @("352.68", "  +5.52%  ")
    6  = @("1.008", "  +0.44%  ")
    7  = @("0.5283", "  +0.99%  ")
    8  = @("0.4565", "  +0.23%  ")
    9  = @("54.32", "  +2.01%  ")
    10 = @("0.09120", "  +2.08%  ")
    11 = @("1.186", "  +0.71%  ")
    12 = @("24.70", "  +2.18%  ")
    13 = @("2.138.61", "  +1.62%  ")
    14 = @("6.907", "  +0.98%  ")
    15 = @("8.149", "  +1.13%  ")
    16 = @("102.28", "  +5.78%  ")
    17 = @("0.00001178", $null)
    18 = @("1.009", "  +0.40%  ")
    19 = @("0.06728", "  +1.12%  ")
    20 = @("19.56", "  +1.63%  ")
    21 = @("1.008", "  +0.44%  ")
    22 = @("6.388", "  +0.66%  ")
    23 = @("30.809.38", "  +1.02%  ")
    24 = @("12.91", "  +3.31%  ")
    25 = @("2.389", "  +1.41%  ")
    26 = @("2.383.19", "  +1.27%  ")
    27 = @("22.60", "  +1.51%  ")
    28 = @("2.582", "  +2.03%  ")
    29 = @("165.13", "  +1.37%  ")
    30 = @("137.23", "  +2.99%  ")
    31 = @("1.213", "  +0.03%  ")
    32 = @("0.1087", "  +1.28%  ")
    33 = @("1.689", "  +0.95%  ")
    34 = @("6.419", "  +0.43%  ")
    35 = @("4.009", "  +1.77%  ")
    36 = @("6.142", "  +7.06%  ")
    37 = @("10.33", "  -0.86%  ")
    38 = @("0.02652", "  +2.35%  ")
    39 = @("0.06939", "  +1.17%  ")
    40 = @("0.2333", "  +1.34%  ")
    41 = @("12.64", "  -0.82%  ")
    42 = @("0.6956", "  +0.95%  ")
    43 = @("1.278", "  +2.21%  ")
    44 = @("14.89", "  +5.96%  ")
    45 = @("2.352", "  +1.59%  ")
    46 = @("0.6479", "  +1.54%  ")
    47 = @($null, "  +2.74%  ")
    48 = @("0.00000000368", "  +5.43%  ")
    49 = @("1.259", "  +0.61%  ")
    50 = @("0.3419", "  +2.15%  ")
    51 = @("83.47", "  +0.00%  ")
}

function Set-TextValue($range, $text) {
    # Force the cell to be stored as text (matching the workbook's original
    # inline-string cells) rather than letting Excel auto-detect numbers,
    # then restore the cell's original (default) style afterwards.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $price = $vals[0]
    $volume = $vals[1]
    if ($null -ne $price) {
        Set-TextValue $ws.Range("D$row") $price
    }
    if ($null -ne $volume) {
        Set-TextValue $ws.Range("E$row") $volume
    }
}
